$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 8: new "Sorting " topic header row (mirrors row 3's style, but across A:C) ----
$ws.Range("A3").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A8").Value = "Sorting "

# ---- Row 9: Selection sort ----
$ws.Range("A4").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").Value = "Selection sort "

$ws.Range("C4").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "goal is with each iteration put smaller element in correct position. example for 5 4 3 2 solution will be 2 4 3 5 ,  2 3 4 5 Remember - find the minimal and swap"

$ws.Range("B9").Value = "https://www.codingninjas.com/studio/problems/selection-sort_624469?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf"
$ws.Hyperlinks.Add($ws.Range("B9"), "https://www.codingninjas.com/studio/problems/selection-sort_624469?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf")
$ws.Range("B4").Copy()
$ws.Range("B9").PasteSpecial(-4122)

# ---- Row 10: Bubble Sort ----
$ws.Range("A5").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = "Bubble Sort"

$ws.Range("B10").Value = "https://www.codingninjas.com/studio/problems/bubble-sort_624380?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf"
$ws.Hyperlinks.Add($ws.Range("B10"), "https://www.codingninjas.com/studio/problems/bubble-sort_624380?utm_source=striver&utm_medium=website&utm_campaign=a_zcoursetuf")
$ws.Range("B5").Copy()
$ws.Range("B10").PasteSpecial(-4122)

# ---- Row heights to match the wrapped multi-line content (rows 4/5 use 75pt) ----
$ws.Rows.Item(9).RowHeight = 75
$ws.Rows.Item(10).RowHeight = 75

# ---- View: scroll down a bit and leave selection on B10 ----
$win = $wb.Windows.Item(1)
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("B10").Select()
